# Fix POC pCO2 extraction mismatch
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Gas any
$ws.Range("C2").Value = 1667
$ws.Range("D2").Value = 316
$ws.Range("E2").Value = 25992
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.06026971329404534
$ws.Range("H2").Value = 0.05752564217369176
$ws.Range("I2").Value = 0.06313591256172663
$ws.Range("J2").Value = 0.8406454866364095
$ws.Range("K2").Value = 0.8238797237297052
$ws.Range("L2").Value = 0.8560940074415263

# Row 3 - ABG threshold
$ws.Range("C3").Value = 830
$ws.Range("D3").Value = 1153
$ws.Range("E3").Value = 3830
$ws.Range("F3").Value = 22162
$ws.Range("G3").Value = 0.1781115879828326
$ws.Range("H3").Value = 0.1673928500551068
$ws.Range("I3").Value = 0.189360584535574
$ws.Range("J3").Value = 0.4185577407967725
$ws.Range("K3").Value = 0.3970227400067738
$ws.Range("L3").Value = 0.4404076706754557

# Row 4 - VBG threshold
$ws.Range("C4").Value = 1465
$ws.Range("D4").Value = 518
$ws.Range("E4").Value = 15441
$ws.Range("F4").Value = 10551
$ws.Range("G4").Value = 0.08665562522181475
$ws.Range("H4").Value = 0.08250821170452956
$ws.Range("I4").Value = 0.09099084007341084
$ws.Range("J4").Value = 0.7387796268280383
$ws.Range("K4").Value = 0.7189959855424114
$ws.Range("L4").Value = 0.7576399311260981

# Row 5 - PCO2 OTHER threshold
$ws.Range("C5").Value = 1229
$ws.Range("D5").Value = 754
$ws.Range("E5").Value = 17280
$ws.Range("F5").Value = 8712
$ws.Range("G5").Value = 0.06640012966664866
$ws.Range("H5").Value = 0.06290242831153363
$ws.Range("I5").Value = 0.07007777703926173
$ws.Range("J5").Value = 0.6197680282400403
$ws.Range("K5").Value = 0.5981897532937894
$ws.Range("L5").Value = 0.6408831721735198
